$d = $word.ActiveDocument

# --- Locate the paragraph index of the anchor paragraph ("Pipes (Pure and
# Impure)") after which the new "Queries" sub-tree must be inserted.
$anchorIdx = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "Pipes (Pure and Impure)") {
        $anchorIdx = $i
    }
    $i = $i + 1
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Appends a new "ListParagraph" bullet right after paragraph $curIdx, at list
# level $ilvl (0-based, matching w:ilvl), with text $text. When $withBreak is
# $true a <w:lastRenderedPageBreak/> marker is written just before the text
# run (mirrors what Word itself records when a page boundary falls there).
# Returns the index of the newly created paragraph.
function Add-Bullet($curIdx, $ilvl, $text, $withBreak) {
    $prevRange = $d.Paragraphs.Item($curIdx).Range
    $prevRange.InsertParagraphAfter()
    $newIdx = $curIdx + 1
    $newRange = $d.Paragraphs.Item($newIdx).Range

    $run = "<w:r><w:t>$text</w:t></w:r>"
    if ($withBreak) {
        $run = "<w:r><w:lastRenderedPageBreak/><w:t>$text</w:t></w:r>"
    }
    $xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"$ilvl`"/><w:numId w:val=`"10`"/></w:numPr></w:pPr>$run</w:p>"
    $newRange.InsertXML($xml)

    return $newIdx
}

# Build the six new bullet paragraphs, in order, right after "Pipes (Pure and
# Impure)" and before "Injection Hierarchy and DI".
$cur = $anchorIdx
$cur = Add-Bullet $cur 1 "Queries" $false
$cur = Add-Bullet $cur 2 "Template variables" $false
$cur = Add-Bullet $cur 2 "ViewChild" $false
$cur = Add-Bullet $cur 2 "ViewChildren" $false
$cur = Add-Bullet $cur 2 "ContentChild" $false
$cur = Add-Bullet $cur 2 "ContentChildren" $true

# The page break that used to render right before "Injection Modifiers" now
# falls earlier in the (now longer) list, inside "ContentChildren" above, so
# strip it from "Injection Modifiers". Re-scan for it (rather than trusting
# the index captured before the insertions above) since the six new bullets
# shifted every later paragraph's index down the document.
$fixIdx = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Injection Modifiers") {
        $fixIdx = $i
    }
    $i = $i + 1
}
$d.Paragraphs.Item($fixIdx).Range.Text = "Injection Modifiers"
